# Adds a new "2022-Q4" quarter sheet (with its own fund-holding table) right
# after the "总计" summary sheet, and records it as a new row on "总计".
#
# Strategy:
#  1. Duplicate the existing "2022-Q3" sheet (tab position 2) to obtain an
#     identically-formatted worksheet inserted right before it (new tab
#     position 2). Rename the duplicate to "2022-Q4".
#  2. Overwrite the duplicate's table with the 2022-Q4 fund data (15 funds,
#     3 more rows than the 12-fund template it was cloned from) and copy the
#     existing row-styling down onto the extra rows.
#  3. On "总计", shift the existing 8 data rows down by one row and write a
#     new first data row for 2022-Q4 (15 holdings, 4.87 亿元).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Step 1: clone "2022-Q3" into a new sheet inserted right before it, then
# rename the clone to "2022-Q4". The clone keeps all of the source sheet's
# formatting (column widths, header styles, borders, etc.) automatically.
# ---------------------------------------------------------------------
$q3 = $wb.Worksheets.Item("2022-Q3")
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item("2022-Q3 (2)")
$q4.Name = "2022-Q4"

# ---------------------------------------------------------------------
# Step 2: replace the cloned table's data with the 2022-Q4 fund holdings.
# Columns: A idx(0-based), B code, C name, D scale, E position, F ratio,
# G value, H rank. Header row (row 1) is left untouched (already correct:
# 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
# ---------------------------------------------------------------------
$q4data = @(
    @("000979", "景顺长城沪港深精选股票",                     "25.01", "92.37", "8.33", "2.0833", 3),
    @("008850", "景顺长城价值稳进三年定期开放灵活配置混合",   "18.18", "97.45", "5.07", "0.9217", 7),
    @("008715", "景顺长城价值驱动一年持有期灵活配置混合",     "6.75",  "92.52", "6.20", "0.4185", 5),
    @("008060", "景顺长城价值边际灵活配置混合A",               "6.42",  "91.73", "6.46", "0.4147", 7),
    @("009098", "景顺长城价值领航两年持有期混合",               "7.08",  "92.98", "5.17", "0.3660", 7),
    @("007291", "汇丰晋信港股通双核策略混合",                   "7.74",  "90.21", "3.91", "0.3026", 5),
    @("015779", "景顺长城价值边际灵活配置混合C",               "2.40",  "91.73", "6.46", "0.1550", 7),
    @("002332", "汇丰晋信沪港深股票A",                           "4.46",  "90.48", "2.52", "0.1124", 9),
    @("002333", "汇丰晋信沪港深股票C",                           "1.79",  "90.48", "2.52", "0.0451", 9),
    @("501310", "华宝标普沪港深中国增强价值指数（LOF）A",       "0.96",  "94.81", "3.01", "0.0289", 7),
    @("004532", "民生加银中证港股通高股息精选指数A",           "0.14",  "92.86", "3.99", "0.0056", 8),
    @("011647", "博时港股通红利精选混合A",                       "0.11",  "77.58", "4.57", "0.0050", 4),
    @("004533", "民生加银中证港股通高股息精选指数C",           "0.09",  "92.86", "3.99", "0.0036", 8),
    @("007397", "华宝标普沪港深中国增强价值指数（LOF）C",       "0.07",  "94.81", "3.01", "0.0021", 7),
    @("011648", "博时港股通红利精选混合C",                       "0.03",  "77.58", "4.57", "0.0014", 4)
)

$rowCount = $q4data.Count
for ($i = 0; $i -lt $rowCount; $i++) {
    $r = $i + 2
    $entry = $q4data[$i]

    $q4.Range("A$r").Value2 = $i
    $q4.Range("B$r").Value2 = $entry[0]
    $q4.Range("C$r").Value2 = $entry[1]
    $q4.Range("D$r").Value2 = $entry[2]
    $q4.Range("E$r").Value2 = $entry[3]
    $q4.Range("F$r").Value2 = $entry[4]
    $q4.Range("G$r").Value2 = $entry[5]
    $q4.Range("H$r").Value2 = $entry[6]
}

# The cloned sheet only had 12 data rows (rows 2-13); stamp the same A-column
# formatting used throughout the table onto the 3 newly-extended rows 14-16.
$q4.Range("A13").Copy()
$q4.Range("A14:A16").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# Step 3: insert the 2022-Q4 row at the top of the "总计" summary table,
# shifting the existing 8 rows down by one (preserving their values, just
# moved + renumbered in column A).
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$existing = @()
for ($r = 2; $r -le 9; $r++) {
    $existing += , @($total.Range("B$r").Value2, $total.Range("C$r").Value2, $total.Range("D$r").Value2)
}

for ($i = $existing.Count - 1; $i -ge 0; $i--) {
    $destRow = $i + 3
    $total.Range("A$destRow").Value2 = $i + 1
    $total.Range("B$destRow").Value2 = $existing[$i][0]
    $total.Range("C$destRow").Value2 = $existing[$i][1]
    $total.Range("D$destRow").Value2 = $existing[$i][2]
}

$total.Range("A2").Value2 = 0
$total.Range("B2").Value2 = "2022-Q4"
$total.Range("C2").Value2 = 15
$total.Range("D2").Value2 = 4.87
